$wb = $excel.ActiveWorkbook

# --- Sheet "Primary Details": insert a new "ID" row right after the header row ---
$ws1 = $wb.Worksheets.Item("Primary Details")

# Shift existing rows 2-6 down to make room for the new row 2 ("ID")
$ws1.Rows.Item(2).Insert()

$ws1.Range("A2").Value = "ID"
$ws1.Range("B2").Value = "66672ce348ccd80bebd64965"

# --- Sheet "Family Information": edit row 2 values ---
$ws5 = $wb.Worksheets.Item("Family Information")

$ws5.Range("A2").Value = "ffedited"
$ws5.Range("B2").Value = "mm-edited"
$ws5.Range("M2").Value = ""
